$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.706.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.582.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.86%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.565'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.589.12'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.038.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.666.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.590.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  +1.15%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0718'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("E31").Value = '  -4.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.74'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.48%  '
$ws.Range("E36").Value = '  -1.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.51'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.30%  '
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.827'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.814'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.50'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '274.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0520'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.47'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.980.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.91%  '
